$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = '44. (Passage Reading: You may think that the leap from naked-eye observation to instrument-aided vision could be one of the greatest advances in the history of the planet. But nobody set out to invent a telescope or eyeglasses. One of the deepest and most widespread human conceptions was faith in the unaided human senses, particularly the naked eye. We do not know who invented eyeglasses, how or where. Everything we know suggests that they were invented by chance, and by a layman who had not been taught about optics. Perhaps an elderly glassmaker making glass disks for leaded windows tested a disk by looking through it and found to his delight that he could see much better. We can suspect that the inventor was not an academic, as professors would delight in boasting of their invention, and before the 13th century, we have no record by any such self-styled inventor. The Italian word lente (English ‘lens’ from ‘lentil’ - the edible seed) or lente di vetro (glass lentil) first used to describe the invention was obviously unacademic; it is clearly not the sort of word that a professor would use to describe the application of his optical theories. From the first recorded uses of eyeglasses before 1300 until the invention of the telescope nearly three hundred years later, lenses were ignored by scholars.) Which of the following could be true of the inventor of eyeglasses?'
$ws.Range("A46").Value = '45. (Passage Reading: You may think that the leap from naked-eye observation to instrument-aided vision could be one of the greatest advances in the history of the planet. But nobody set out to invent a telescope or eyeglasses. One of the deepest and most widespread human conceptions was faith in the unaided human senses, particularly the naked eye. We do not know who invented eyeglasses, how or where. Everything we know suggests that they were invented by chance, and by a layman who had not been taught about optics. Perhaps an elderly glassmaker making glass disks for leaded windows tested a disk by looking through it and found to his delight that he could see much better. We can suspect that the inventor was not an academic, as professors would delight in boasting of their invention, and before the 13th century, we have no record by any such self-styled inventor. The Italian word lente (English ‘lens’ from ‘lentil’ - the edible seed) or lente di vetro (glass lentil) first used to describe the invention was obviously unacademic; it is clearly not the sort of word that a professor would use to describe the application of his optical theories. From the first recorded uses of eyeglasses before 1300 until the invention of the telescope nearly three hundred years later, lenses were ignored by scholars.) What makes the author think that the invention of eyeglasses was unacademic?'
$ws.Range("A47").Value = '46. (Passage Reading: You may think that the leap from naked-eye observation to instrument-aided vision could be one of the greatest advances in the history of the planet. But nobody set out to invent a telescope or eyeglasses. One of the deepest and most widespread human conceptions was faith in the unaided human senses, particularly the naked eye. We do not know who invented eyeglasses, how or where. Everything we know suggests that they were invented by chance, and by a layman who had not been taught about optics. Perhaps an elderly glassmaker making glass disks for leaded windows tested a disk by looking through it and found to his delight that he could see much better. We can suspect that the inventor was not an academic, as professors would delight in boasting of their invention, and before the 13th century, we have no record by any such self-styled inventor. The Italian word lente (English ‘lens’ from ‘lentil’ - the edible seed) or lente di vetro (glass lentil) first used to describe the invention was obviously unacademic; it is clearly not the sort of word that a professor would use to describe the application of his optical theories. From the first recorded uses of eyeglasses before 1300 until the invention of the telescope nearly three hundred years later, lenses were ignored by scholars.) What is the main idea of the passage?'
$ws.Range("A49").Value = '48.  (Passage Reading: The annual melt of Arctic ice is under way, and the melt is getting worse and worse. The size of the ice cap where the melt is taking place has become the smallest since satellite records began. It will not be many more years until the Arctic is ice-free in summer. This is not news. Nor is the fact that the US is about to break its promise on its climate change policies, again. The scale and speed of worrying events in the Arctic is increasing, so is the recklessness of the current administration’s climate policies. The two are closely connected: without comprehensive, rapid, global action, there is zero chance of preventing irreversible change in the Arctic, which may bring catastrophe to the rest of the world. What is also not-news is the antibiotic crisis. For years, if not decades, medics have been warning that overuse is driving up resistance, threating to propel us back to a pre-penicillin age when common bacterial infections could be fatal. Both are large, complex problems arising from technological progress. Both are tragedies that require collective action. But on antibiotic resistance, there are signs that the problem is getting less severe, thanks largely to changing prescription habits, following great collective action, which is missing in the climate change problem.) What is the author’s attitude towards the current US administration’s climate policies?'
$ws.Range("A50").Value = '49.  (Passage Reading: The annual melt of Arctic ice is under way, and the melt is getting worse and worse. The size of the ice cap where the melt is taking place has become the smallest since satellite records began. It will not be many more years until the Arctic is ice-free in summer. This is not news. Nor is the fact that the US is about to break its promise on its climate change policies, again. The scale and speed of worrying events in the Arctic is increasing, so is the recklessness of the current administration’s climate policies. The two are closely connected: without comprehensive, rapid, global action, there is zero chance of preventing irreversible change in the Arctic, which may bring catastrophe to the rest of the world. What is also not-news is the antibiotic crisis. For years, if not decades, medics have been warning that overuse is driving up resistance, threating to propel us back to a pre-penicillin age when common bacterial infections could be fatal. Both are large, complex problems arising from technological progress. Both are tragedies that require collective action. But on antibiotic resistance, there are signs that the problem is getting less severe, thanks largely to changing prescription habits, following great collective action, which is missing in the climate change problem.) It is clearly stated in the passage that the antibiotic crisis ----.'
$ws.Range("A51").Value = '50.  (Passage Reading: The annual melt of Arctic ice is under way, and the melt is getting worse and worse. The size of the ice cap where the melt is taking place has become the smallest since satellite records began. It will not be many more years until the Arctic is ice-free in summer. This is not news. Nor is the fact that the US is about to break its promise on its climate change policies, again. The scale and speed of worrying events in the Arctic is increasing, so is the recklessness of the current administration’s climate policies. The two are closely connected: without comprehensive, rapid, global action, there is zero chance of preventing irreversible change in the Arctic, which may bring catastrophe to the rest of the world. What is also not-news is the antibiotic crisis. For years, if not decades, medics have been warning that overuse is driving up resistance, threating to propel us back to a pre-penicillin age when common bacterial infections could be fatal. Both are large, complex problems arising from technological progress. Both are tragedies that require collective action. But on antibiotic resistance, there are signs that the problem is getting less severe, thanks largely to changing prescription habits, following great collective action, which is missing in the climate change problem.) What is the main aim of the author in writing this passage?'
$ws.Range("A53").Value = '52. (Passage Reading: Child labour has traditionally referred to the use of children to perform adult work on farms and in factories and mines. Ironically, it has helped as well as hindered the development of universal public education in the United States. Historically, this problem did not become a controversial social issue until the beginning of the Industrial Revolution in the late 18th century. After a series of wars produced a shortage of adult male labour in Britain, cotton mill owners offered children free room and meals in company-owned dormitories in exchange for their labour. Government authorities were delighted to get homeless children off the streets. For the mills, the introduction of child labour proved one of the most profitable innovations of the Industrial Revolution. The practice spread to the United States, which also faced a shortage of adult male factory workers to produce food and tobacco, the staples of the new nation’s trade. At the time, most schools were already reserved for the children of landowners, and therefore, an entire generation of child labourers had grown up illiterate by the 1830s. As the quality of goods they produced deteriorated and overseas sales of American goods declined, many of the very industrialists who had exploited child labour at the turn of the century joined social reformers who wanted to end such exploitation, supporting universal compulsory education to ensure that their future workers would be able to read, write and calculate.) What caused child labourers in the United States to become illiterate?'
$ws.Range("A54").Value = '53. (Passage Reading: Child labour has traditionally referred to the use of children to perform adult work on farms and in factories and mines. Ironically, it has helped as well as hindered the development of universal public education in the United States. Historically, this problem did not become a controversial social issue until the beginning of the Industrial Revolution in the late 18th century. After a series of wars produced a shortage of adult male labour in Britain, cotton mill owners offered children free room and meals in company-owned dormitories in exchange for their labour. Government authorities were delighted to get homeless children off the streets. For the mills, the introduction of child labour proved one of the most profitable innovations of the Industrial Revolution. The practice spread to the United States, which also faced a shortage of adult male factory workers to produce food and tobacco, the staples of the new nation’s trade. At the time, most schools were already reserved for the children of landowners, and therefore, an entire generation of child labourers had grown up illiterate by the 1830s. As the quality of goods they produced deteriorated and overseas sales of American goods declined, many of the very industrialists who had exploited child labour at the turn of the century joined social reformers who wanted to end such exploitation, supporting universal compulsory education to ensure that their future workers would be able to read, write and calculate.) What was the real reason why industrialists started to promote universal compulsory education?'
$ws.Range("A55").Value = '54. (Passage Reading: Child labour has traditionally referred to the use of children to perform adult work on farms and in factories and mines. Ironically, it has helped as well as hindered the development of universal public education in the United States. Historically, this problem did not become a controversial social issue until the beginning of the Industrial Revolution in the late 18th century. After a series of wars produced a shortage of adult male labour in Britain, cotton mill owners offered children free room and meals in company-owned dormitories in exchange for their labour. Government authorities were delighted to get homeless children off the streets. For the mills, the introduction of child labour proved one of the most profitable innovations of the Industrial Revolution. The practice spread to the United States, which also faced a shortage of adult male factory workers to produce food and tobacco, the staples of the new nation’s trade. At the time, most schools were already reserved for the children of landowners, and therefore, an entire generation of child labourers had grown up illiterate by the 1830s. As the quality of goods they produced deteriorated and overseas sales of American goods declined, many of the very industrialists who had exploited child labour at the turn of the century joined social reformers who wanted to end such exploitation, supporting universal compulsory education to ensure that their future workers would be able to read, write and calculate.) What is the primary purpose of the author?'
$ws.Range("A57").Value = '56. (Passage Reading: Earth is teeming with life — and water makes it all possible. But elsewhere in the cosmos, life might be built from different chemicals that dissolve and assemble in some other liquid: perhaps methane, kerosene, or even chloroform. For now, it is not feasible for humans to set foot on those worlds and see what is there, but researchers are exploring some possibilities in labs here on Earth. The idea of a strange parallel form of life, whose cells do the same basic things as ours using a completely different chemistry, is not new to science. In fact, Isaac Asimov broached the subject in his 1962 essay “Not as We Know It: The Chemistry of Life.” And in 2004, the same year the Cassini spacecraft entered Saturn’s orbit, biochemist Steven Benner proposed that on a world like Saturn’s moon Titan, life might use liquid hydrocarbons as a solvent (a liquid that can dissolve other substances), the way water is used on Earth. With new exoplanets joining the roster of known worlds every day, it is likely that some of them have oceans (or at least warm puddles) of hexane, ethers, chloroform, or other exotic liquids that might serve as a source of life, which we do not yet know.) It is stated in the passage that Isaac Asimov ----.'
$ws.Range("A58").Value = '57. (Passage Reading: Earth is teeming with life — and water makes it all possible. But elsewhere in the cosmos, life might be built from different chemicals that dissolve and assemble in some other liquid: perhaps methane, kerosene, or even chloroform. For now, it is not feasible for humans to set foot on those worlds and see what is there, but researchers are exploring some possibilities in labs here on Earth. The idea of a strange parallel form of life, whose cells do the same basic things as ours using a completely different chemistry, is not new to science. In fact, Isaac Asimov broached the subject in his 1962 essay “Not as We Know It: The Chemistry of Life.” And in 2004, the same year the Cassini spacecraft entered Saturn’s orbit, biochemist Steven Benner proposed that on a world like Saturn’s moon Titan, life might use liquid hydrocarbons as a solvent (a liquid that can dissolve other substances), the way water is used on Earth. With new exoplanets joining the roster of known worlds every day, it is likely that some of them have oceans (or at least warm puddles) of hexane, ethers, chloroform, or other exotic liquids that might serve as a source of life, which we do not yet know.) The underlined word in the passage ‘broached’ is closest in meaning to ----.'
$ws.Range("A59").Value = '58. (Passage Reading: Earth is teeming with life — and water makes it all possible. But elsewhere in the cosmos, life might be built from different chemicals that dissolve and assemble in some other liquid: perhaps methane, kerosene, or even chloroform. For now, it is not feasible for humans to set foot on those worlds and see what is there, but researchers are exploring some possibilities in labs here on Earth. The idea of a strange parallel form of life, whose cells do the same basic things as ours using a completely different chemistry, is not new to science. In fact, Isaac Asimov broached the subject in his 1962 essay “Not as We Know It: The Chemistry of Life.” And in 2004, the same year the Cassini spacecraft entered Saturn’s orbit, biochemist Steven Benner proposed that on a world like Saturn’s moon Titan, life might use liquid hydrocarbons as a solvent (a liquid that can dissolve other substances), the way water is used on Earth. With new exoplanets joining the roster of known worlds every day, it is likely that some of them have oceans (or at least warm puddles) of hexane, ethers, chloroform, or other exotic liquids that might serve as a source of life, which we do not yet know.) What is the passage mainly about?'
$ws.Range("A61").Value = '60. (Passage Reading: Companies often link pay increases to a measure of performance to motivate and retain the best employees, which may also ensure the companies’ long-term success. Among many forms of individual performance-related pay, the most popular is merit-based pay. Under such a system, an employee’s pay increase is based on a supervisor’s assessment of performance over a previous time period, usually the prior 12 months. The annual cycle of the pay scheme begins with the establishment of performance objectives, either by the supervisor or jointly between the employee and the supervisor. Progress is monitored over the ensuing months, and the supervisor rates the performance of an employee, on a given scale, against these objectives. Each of these rating points is then associated with the amount of a pay increase, expressed as a percentage, so the higher the rating, the higher the pay increase. The size of the pay increase can also vary depending on the position in the pay range. Typically, employees at the bottom of the pay range with high-performance ratings will get a larger increase than an employee near the top of the pay range with a high-performance rating. In some cases, the intention is to enable the employee to reach the top of the pay range at the same time as he or she is ready to be promoted. The difficulty is that this system is not supported by all employees. Research has demonstrated that merit-based pay is attractive to some occupations, for example managers, and to certain demographic groups, such as white-collar workers, high-income earners, and younger employees, but it might not be the best payment approach for all occupations and companies.) It is understood from the passage that in the merit-based pay system, ----.'
$ws.Range("A63").Value = '62. (Passage Reading: When we take medication, some can pass through us intact and go down the toilet. Once in the sewers, these drugs are not removed during the treatment process, and end up in the waterways. To find out how pharmaceutical waste finds its way into aquatic creatures, Erinn Richmond at Monash University in Australia and her colleagues sampled flies, beetles, other insects and spiders from six waterways in the greater Melbourne region. Overall, the group detected 69 medications in the invertebrates. Organisms near wastewater plants contained the highest levels, but low levels were also detected in those from more pristine areas. Richmond presumes river-borne pharmaceuticals probably accumulate in flies and beetles while they are underwater larvae, then transfer to spiders that feed on them after they emerge as adults. Other predators in the food web such as fish, platypuses, birds, bats and frogs may also become cross-contaminated, she says. Her team calculated that trout and platypuses feeding on insects near wastewater plants would consume anti-depressant doses up to half the level normally prescribed to humans. The impacts of this are not certain, but previous studies have found that high-level exposure to anti-depressants causes fish to become less vigilant about avoiding predators. Male birds also seem to sing less to females when exposed to the drugs.) What could be the best title for this passage?'
